$d = $word.ActiveDocument

# Helper: run a single-occurrence Find/Replace confined to a [start,end)
# character range. We re-derive a fresh Document.Range(start,end) for the
# actual Find call (rather than calling .Find on the Cell/Paragraph Range
# object directly) because Find.Execute on those derived Range objects
# ignores their bounds and searches from the top of the story.
function Replace-InRange {
    param(
        [int]$StartPos,
        [int]$EndPos,
        [string]$FindText,
        [string]$ReplaceText,
        [bool]$ReplaceAll = $false
    )
    $r = $d.Range($StartPos, $EndPos)
    $mode = 1
    if ($ReplaceAll) { $mode = 2 }
    return $r.Find.Execute($FindText, $true, $false, $false, $false, $false, `
                            $true, 0, $false, $ReplaceText, $mode)
}

function Replace-InCell {
    param(
        [__ComObject]$Table,
        [int]$Row,
        [int]$Col,
        [string]$FindText,
        [string]$ReplaceText
    )
    $c = $Table.Cell($Row, $Col)
    Replace-InRange $c.Range.Start $c.Range.End $FindText $ReplaceText $false | Out-Null
}

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$FindText,
        [string]$ReplaceText
    )
    $p = $d.Paragraphs.Item($Index)
    Replace-InRange $p.Range.Start $p.Range.End $FindText $ReplaceText $false | Out-Null
}

# 1) Digital seal value (unique in document)
$d.Content.Find.Execute("51679312_", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "86199828_", 2) | Out-Null

# 2) Paragraph 8: "... НЕО G6 заводской № 10250112 соответствует ..."
Replace-InParagraph 8 " G6" " G25"
Replace-InParagraph 8 "10250112 " "10250324 "

# 3) Paragraph 26: "...счетчика НЕО G6          №10250112   в соответствии..."
Replace-InParagraph 26 " G6" " G25"
Replace-InParagraph 26 "10250112" "10250324"

# 4) Paragraph 37: bold heading "НЕО G6" right before the spec table
Replace-InParagraph 37 "G6" "G25"

$t = $d.Tables.Item(2)

# 5) Row 1: table header "НЕО G6"
Replace-InCell $t 1 2 "G6" "G25"

# 6) Row 2: model code "НЕО-G6-1-*-А-*-2-Л-250-*" -> "НЕО-G25-2-P-А-*-16-Л-335-*"
Replace-InCell $t 2 2 "G6" "G25"
Replace-InCell $t 2 2 "-1-" "-2-"
Replace-InCell $t 2 2 "*" "P"
Replace-InCell $t 2 2 "2-Л" "16-Л"
Replace-InCell $t 2 2 "250" "335"

# 7) Remaining characteristic rows
Replace-InCell $t 3 2 "0,06" "0,25"
Replace-InCell $t 4 2 "6" "25"
Replace-InCell $t 5 2 "10" "40"
Replace-InCell $t 6 2 "0,008" "0,01"
Replace-InCell $t 7 2 "250" "300"
Replace-InCell $t 10 2 "2" "16"
Replace-InCell $t 16 2 "335" "468"
Replace-InCell $t 17 2 "243" "430"
Replace-InCell $t 18 2 "182" "289"
Replace-InCell $t 19 2 "250" "335"
Replace-InCell $t 20 2 "1 1/4" "2 1/2"
Replace-InCell $t 21 2 "3,9" "10,6"

Write-Host "Done."
